$d = $word.ActiveDocument

# 1. Update the date stamp in the document header block.
$d.Content.Find.Execute("2024-06-21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-11-28", 2)

# 2. The code-availability sentence ends "...and code is available from
#    [the OSF]" where "the OSF" is a hyperlink whose URL is hidden (not
#    visible as literal text). Make the OSF repo URL visible as plain text
#    instead of hiding it behind link text, and drop the hyperlink field.
#    Result: "...and code is available from the OSF at https://osf.io/629bx."
$target = $null
foreach ($h in $d.Hyperlinks) {
    if ($h.TextToDisplay -eq "the OSF") {
        $target = $h
    }
}

if ($target -ne $null) {
    $osfUrl = $target.Address
    $target.Delete()

    $searchRange = $d.Content
    $f = $searchRange.Find
    $f.ClearFormatting()
    $f.Text = "and code is available from the OSF."
    $found = $f.Execute()
    if ($found) {
        $searchRange.Delete()
        $insertion = $searchRange.Duplicate
        $insertion.Collapse(1)
        $insertion.InsertAfter("and code is available from the OSF at $osfUrl.")
    }
}
